$d = $word.ActiveDocument

# --- Change 1: "Models" paragraph -> "Models / Layers" + new explanatory paragraphs ---
$rng1 = $d.Content
$found = $rng1.Find.Execute("Models")
if (-not $found) { throw "Could not find 'Models' anchor text" }
$p1a = $rng1.Paragraphs(1)
$p1b = $p1a.Next()
$replaceRange1 = $d.Range($p1a.Range.Start, $p1b.Range.End)
$xml1 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="start"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Models / </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Layers</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="start"/><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="start"/><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>The idea is to enable model representations being equivalent (containing the same data) in various layers to be switched back an forth between each layer representation to be used in the most appropriate task for a given representation.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="start"/><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$replaceRange1.InsertXML($xml1)

# --- Change 2: "- id : ID" (inside the IDOccurrence section) -> "- occurringId : ID" (split runs) ---
$rng2 = $d.Content
$found = $rng2.Find.Execute("IDOccurrence : ID")
if (-not $found) { throw "Could not find 'IDOccurrence : ID' anchor text" }
$rng2.Collapse(0)
$found = $rng2.Find.Execute("- id : ID")
if (-not $found) { throw "Could not find '- id : ID' following IDOccurrence heading" }
$p2 = $rng2.Paragraphs(1)
$replaceRange2 = $d.Range($p2.Range.Start, $p2.Range.End)
$xml2 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="start"/><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>occurringI</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>d : ID</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$replaceRange2.InsertXML($xml2)

# --- Change 3: merge "Each ID assigned..." + "Object Context..." paragraphs, expand text, add new paragraphs ---
$rng3 = $d.Content
$found = $rng3.Find.Execute("Each ID assigned")
if (-not $found) { throw "Could not find 'Each ID assigned' anchor text" }
$p3a = $rng3.Paragraphs(1)
$p3b = $p3a.Next()
$replaceRange3 = $d.Range($p3a.Range.Start, $p3b.Range.End)
$xml3 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="start"/><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Each ID </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">is </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">assigned a unique prime number ID </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>at creation time</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">. FCA Context / Lattices built upon, for example for a given Data / Schema predicate / arc occurrence role, having the context objects being the statement occurrence subjects and the context attributes the statement occurrence objects, Predicate </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>FCA</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> Context: (Subjects x Objects). For a subject statement occurrence the context is: Subject </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>FCA</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> Context: (Predicates x Objects and for an object statement occurrence role the context is: Object </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>FCA</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> Context (Subjectx x Predicates).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="start"/><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="start"/><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Embeddings: For an ID, its prime ID number plus all ID’s occurrences embeddings. For an IDOccurrence, its ID class embeddings, its occurring ID embeddings and its context embeddings. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="start"/><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="start"/><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="arial" w:hAnsi="arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Embeddings similarity: IDs, IDOccurrences sharing the same primes for their embeddings in a given context. FCA Concept Lattice Clustering. (TODO). </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$replaceRange3.InsertXML($xml3)

Write-Output "Done. Paragraphs.Count=$($d.Paragraphs.Count)"
